# Update the F20 Stat 5100 schedule through module 4.
# The only substantive content change is row 47 (12/11 class day), which
# changes from "7.3: Modern Applications" to "No Class (Interim Day)" and
# is formatted to match the other "No Class (...)" entries (bold italic).
# A few neighboring cells (Final Project Introductions / Model
# Interpretabillity / Group Work Primer) are re-entered with a trailing
# space, matching how the author's update touched those notes as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B47").Value = "No Class (Interim Day)"
$ws.Range("B47").Font.Bold = $true
$ws.Range("B47").Font.Italic = $true

$ws.Range("B31").Value = "Group Work Primer "
$ws.Range("B30").Value = "Model Interpretabillity "
$ws.Range("B29").Value = "Final Project Introductions "

$ws.Range("G32").Select()
